$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 290.57144
$ws.Range("I2").Value = 167.6
$ws.Range("J2").Value = 358.8889
$ws.Range("K2").Value = 167.6
$ws.Range("L2").Value = 358.8889
$ws.Range("M2").Value = -54.59999999999999
$ws.Range("N2").Value = -584.8888999999999
$ws.Range("H41").Value = 1466.7646
$ws.Range("I41").Value = 1668.5714
$ws.Range("J41").Value = 1325.5
$ws.Range("K41").Value = 1668.5714
$ws.Range("L41").Value = 1325.5
$ws.Range("M41").Value = -1228.5714
$ws.Range("N41").Value = -2205.5
$ws.Range("H58").Value = 1034.2727
$ws.Range("J58").Value = 1686.1666
$ws.Range("L58").Value = 5058.4998
$ws.Range("N58").Value = -5358.4998
$ws.Range("H127").Value = 1072.8572
$ws.Range("I127").Value = 502
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 1506
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 3454
$ws.Range("N127").Value = -17420
$ws.Range("H137").Value = 1300.7307
$ws.Range("J137").Value = 2198.625
$ws.Range("L137").Value = 6595.875
$ws.Range("N137").Value = -11695.875
$ws.Range("H138").Value = 1581.1837
$ws.Range("I138").Value = 841.93335
$ws.Range("J138").Value = 1907.3235
$ws.Range("K138").Value = 2525.80005
$ws.Range("L138").Value = 5721.970499999999
$ws.Range("M138").Value = 2614.19995
$ws.Range("N138").Value = -16001.9705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2862.5945
$ws.Range("I32").Value = 2837.5
$ws.Range("K32").Value = 2837.5
$ws.Range("M32").Value = -2550.5
$ws.Range("H45").Value = 1042.3334
$ws.Range("I45").Value = 1025.9333
$ws.Range("J45").Value = 1083.3334
$ws.Range("K45").Value = 1025.9333
$ws.Range("L45").Value = 1083.3334
$ws.Range("M45").Value = -648.9332999999999
$ws.Range("N45").Value = -1837.3334
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 20836762
$ws.Range("I102").Value = 23812870
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 23812870
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -23811248
$ws.Range("N102").Value = -7244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1719.8334
$ws.Range("I20").Value = 1373
$ws.Range("J20").Value = 2066.6667
$ws.Range("K20").Value = 1373
$ws.Range("L20").Value = 2066.6667
$ws.Range("M20").Value = -1126
$ws.Range("N20").Value = -2560.6667
$ws.Range("H94").Value = 16667399
$ws.Range("I94").Value = 22727994
$ws.Range("J94").Value = 762.5
$ws.Range("K94").Value = 22727994
$ws.Range("L94").Value = 762.5
$ws.Range("M94").Value = -22727543
$ws.Range("N94").Value = -1664.5
$ws.Range("H99").Value = 38462770
$ws.Range("I99").Value = 50001000
$ws.Range("J99").Value = 1996.6666
$ws.Range("K99").Value = 50001000
$ws.Range("L99").Value = 1996.6666
$ws.Range("M99").Value = -49999502
$ws.Range("N99").Value = -4992.6666
$ws.Range("H105").Value = 83335920
$ws.Range("I105").Value = 83335920
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 83335920
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -83334173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.5
$ws.Range("J7").Value = 88
$ws.Range("L7").Value = 88
$ws.Range("N7").Value = -314
$ws.Range("H62").Value = 9526133
$ws.Range("I62").Value = 2440
$ws.Range("J62").Value = 200000000
$ws.Range("K62").Value = 2440
$ws.Range("L62").Value = 200000000
$ws.Range("M62").Value = -1816
$ws.Range("N62").Value = -200001248
$ws.Range("H65").Value = 9526133
$ws.Range("I65").Value = 2440
$ws.Range("J65").Value = 200000000
$ws.Range("K65").Value = 12200
$ws.Range("L65").Value = 1000000000
$ws.Range("M65").Value = -9080
$ws.Range("N65").Value = -1000006240
$ws.Range("H99").Value = 1661.8334
$ws.Range("I99").Value = 1749
$ws.Range("J99").Value = 1539.8
$ws.Range("K99").Value = 1749
$ws.Range("L99").Value = 1539.8
$ws.Range("M99").Value = -251
$ws.Range("N99").Value = -4535.8
$ws.Range("H105").Value = 861.8
$ws.Range("I105").Value = 770
$ws.Range("J105").Value = 999.5
$ws.Range("K105").Value = 770
$ws.Range("L105").Value = 999.5
$ws.Range("M105").Value = 977
$ws.Range("N105").Value = -4493.5
$ws.Range("H107").Value = 382.52942
$ws.Range("I107").Value = 238.58333
$ws.Range("K107").Value = 238.58333
$ws.Range("M107").Value = 1681.41667
$ws.Range("H126").Value = 1661.8334
$ws.Range("I126").Value = 1749
$ws.Range("J126").Value = 1539.8
$ws.Range("K126").Value = 5247
$ws.Range("L126").Value = 4619.4
$ws.Range("M126").Value = -2777
$ws.Range("N126").Value = -9559.4
$ws.Range("H141").Value = 27950.75
$ws.Range("J141").Value = 27950.75
$ws.Range("L141").Value = 27950.75
$ws.Range("N141").Value = -38310.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.7
$ws.Range("J2").Value = 124.5
$ws.Range("L2").Value = 747
$ws.Range("N2").Value = -973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62.875
$ws.Range("I2").Value = 13.833333
$ws.Range("J2").Value = 210
$ws.Range("K2").Value = 13.833333
$ws.Range("L2").Value = 210
$ws.Range("M2").Value = 99.166667
$ws.Range("N2").Value = -436
$ws.Range("H80").Value = 2777.7778
$ws.Range("I80").Value = 1533.3334
$ws.Range("J80").Value = 5266.6665
$ws.Range("K80").Value = 1533.3334
$ws.Range("L80").Value = 5266.6665
$ws.Range("M80").Value = -535.3334
$ws.Range("N80").Value = -7262.6665
$ws.Range("H83").Value = 2777.7778
$ws.Range("I83").Value = 1533.3334
$ws.Range("J83").Value = 5266.6665
$ws.Range("K83").Value = 7666.666999999999
$ws.Range("L83").Value = 26333.3325
$ws.Range("M83").Value = -2674.666999999999
$ws.Range("N83").Value = -36317.3325
$ws.Range("H97").Value = 751.6667
$ws.Range("I97").Value = 751.6667
$ws.Range("K97").Value = 751.6667
$ws.Range("M97").Value = -255.6667
$ws.Range("H122").Value = 1460.1818
$ws.Range("I122").Value = 1605.375
$ws.Range("J122").Value = 1073
$ws.Range("K122").Value = 4816.125
$ws.Range("L122").Value = 3219
$ws.Range("M122").Value = -2366.125
$ws.Range("N122").Value = -8119
$ws.Range("H126").Value = 2104.4443
$ws.Range("I126").Value = 1751.4546
$ws.Range("J126").Value = 2659.1428
$ws.Range("K126").Value = 5254.3638
$ws.Range("L126").Value = 7977.428400000001
$ws.Range("M126").Value = -2784.3638
$ws.Range("N126").Value = -12917.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1610.1875
$ws.Range("I7").Value = 1314.3636
$ws.Range("J7").Value = 2261
$ws.Range("K7").Value = 1314.3636
$ws.Range("L7").Value = 2261
$ws.Range("M7").Value = -1202.3636
$ws.Range("N7").Value = -2485
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("N36").Value = 0
$ws.Range("H40").Value = 3199.4546
$ws.Range("I40").Value = 3026
$ws.Range("J40").Value = 3298.5715
$ws.Range("K40").Value = 3026
$ws.Range("L40").Value = 3298.5715
$ws.Range("M40").Value = -2890
$ws.Range("N40").Value = -3570.5715
$ws.Range("H61").Value = 1455.5714
$ws.Range("I61").Value = 1161.8
$ws.Range("K61").Value = 1161.8
$ws.Range("M61").Value = -959.8
$ws.Range("H93").Value = 754.8570999999999
$ws.Range("I93").Value = 754.8570999999999
$ws.Range("K93").Value = 754.8570999999999
$ws.Range("M93").Value = 493.1429000000001
$ws.Range("H100").Value = 2380.6
$ws.Range("I100").Value = 2249
$ws.Range("J100").Value = 2413.5
$ws.Range("K100").Value = 2249
$ws.Range("L100").Value = 2413.5
$ws.Range("M100").Value = -1708
$ws.Range("N100").Value = -3495.5
$ws.Range("H101").Value = 16999.5
$ws.Range("J101").Value = 16999.5
$ws.Range("L101").Value = 16999.5
$ws.Range("N101").Value = -23489.5
$ws.Range("H113").Value = 1455.5714
$ws.Range("I113").Value = 1161.8
$ws.Range("K113").Value = 1161.8
$ws.Range("M113").Value = 1008.2
$ws.Range("H122").Value = 50002180
$ws.Range("I122").Value = 62502124
$ws.Range("J122").Value = 2405
$ws.Range("K122").Value = 187506372
$ws.Range("L122").Value = 7215
$ws.Range("M122").Value = -187503922
$ws.Range("N122").Value = -12115
$ws.Range("H126").Value = 1610.1875
$ws.Range("I126").Value = 1314.3636
$ws.Range("J126").Value = 2261
$ws.Range("K126").Value = 3943.0908
$ws.Range("L126").Value = 6783
$ws.Range("M126").Value = -1473.0908
$ws.Range("N126").Value = -11723
$ws.Range("H132").Value = 30259
$ws.Range("I132").Value = 1062.9
$ws.Range("K132").Value = 3188.7
$ws.Range("M132").Value = -658.7000000000003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 14500
$ws.Range("J76").Value = 14500
$ws.Range("L76").Value = 14500
$ws.Range("N76").Value = -15130
$ws.Range("H79").Value = 14500
$ws.Range("J79").Value = 14500
$ws.Range("L79").Value = 14500
$ws.Range("N79").Value = -16684
$ws.Range("H96").Value = 1860.8572
$ws.Range("I96").Value = 1811.8
$ws.Range("J96").Value = 1983.5
$ws.Range("K96").Value = 1811.8
$ws.Range("L96").Value = 1983.5
$ws.Range("M96").Value = -438.8
$ws.Range("N96").Value = -4729.5
$ws.Range("H113").Value = 886.3333
$ws.Range("I113").Value = 518.2
$ws.Range("J113").Value = 1346.5
$ws.Range("K113").Value = 1554.6
$ws.Range("L113").Value = 4039.5
$ws.Range("M113").Value = 615.3999999999999
$ws.Range("N113").Value = -8379.5
$ws.Range("H122").Value = 12382357
$ws.Range("I122").Value = 15295454
$ws.Range("J122").Value = 1698.75
$ws.Range("K122").Value = 45886362
$ws.Range("L122").Value = 5096.25
$ws.Range("M122").Value = -45883912
$ws.Range("N122").Value = -9996.25
$ws.Range("H126").Value = 76924570
$ws.Range("I126").Value = 76924570
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 230773710
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -230771240
